# Feature Tracker / Install Tracker update
#
# Mirrors the author's edit:
#   - Mark "Roll groupings" (row 5, completed in 1.7.0), "Add Genesys dice"
#     (row 17) and "Add Fantasy Flight Star Wars dice" (row 18) as filtered
#     out (they already had a Completed Version, they just weren't hidden
#     by the autofilter yet in the source file).
#   - Finish two more features - "Quick category change" (row 27) and
#     "Compact view - Saved" (row 28) - by stamping them with Completed
#     Version "1.7.2" in column C.
#   - Because the "only show blanks in Completed Version" autofilter now
#     matches five more completed rows, hide those rows and grow the
#     filtered/named range from row 26 to row 29.
#   - Leave the cursor where the author left it after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Stamp the newly-completed features with their Completed Version
#    *before* hiding their rows - writing a value into a previously empty
#    cell of an already-hidden row makes the host recompute/override the
#    row height (ht/customHeight), which the source file doesn't have.
$ws.Range("C27").Value = "1.7.2"
$ws.Range("C28").Value = "1.7.2"

# 2) Hide all rows that are now "completed" (Completed Version not blank)
#    so they match the existing "show blanks only" autofilter state.
$ws.Rows(5).Hidden = $true
$ws.Rows(17).Hidden = $true
$ws.Rows(18).Hidden = $true
$ws.Rows(27).Hidden = $true
$ws.Rows(28).Hidden = $true

# 3) Grow the hidden _FilterDatabase defined name to the new data extent.
$filterDatabase = $wb.Names.Item("Sheet1!_FilterDatabase")
$filterDatabase.RefersTo = "=Sheet1!`$A`$1:`$E`$29"

# 4) Rebuild the worksheet autofilter over the grown range A1:E29, keeping
#    the same "blanks only" filter on column C (field 3 / colId 2).
$ws.AutoFilterMode = $false
$blankOnly = @("")
$ws.Range("A1:E29").AutoFilter(3, $blankOnly, 7) | Out-Null

# 5) Restore the author's final selection/cursor position.
$ws.Range("B34").Select() | Out-Null
